# Update battery standby / alarm load values for the new accessories test case (NGC1928)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("F8").Value = 0.3
$ws.Range("G8").Value = 0.612

# Reflect the selection/view state saved with the workbook
$ws.Activate() | Out-Null
$ws.Range("F8").Select() | Out-Null
